$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = @((1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02),(1.02))
$rng = $ws.Range("B2:B25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @((1.018245970556423),(1.019142798033991),(1.019723537735224),(1.019967782910709),(1.020008798705983),(1.019726800951245),(1.018548967296559),(1.016476842703709),(1.015097767876342),(1.014501184491556),(1.014279673001967),(1.014327184050356),(1.014482872518264),(1.014578808773562),(1.015137373162919),(1.015487898019444),(1.01569240768002),(1.015762149385926),(1.015450284379481),(1.014437023756973),(1.013800445596852),(1.014137860008781),(1.015467280196209),(1.017012129478322))
$rng = $ws.Range("C2:C25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @((1.023935713125387),(1.024579609787356),(1.024996328756963),(1.025171533601696),(1.02520095216563),(1.024998669791218),(1.024153304507182),(1.022664311004803),(1.021672185919217),(1.021242733556713),(1.021083239224717),(1.021117450214172),(1.021229549218312),(1.02129862026589),(1.021700690445452),(1.021952938300908),(1.022100084020164),(1.022150259198655),(1.021925873052176),(1.021196538196005),(1.020738112635157),(1.02098111921743),(1.021938102630586),(1.023049164237387))
$rng = $ws.Range("D2:D25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @((1.01953932210469),(1.02029831003834),(1.020790200135669),(1.020997174442783),(1.021031937065295),(1.020792965016679),(1.01979566436168),(1.018044293925567),(1.016880844497104),(1.016378058057887),(1.016191451877276),(1.016231472673268),(1.016362630041212),(1.016443460493757),(1.016914233876177),(1.017209805076963),(1.017382302636523),(1.017441136059809),(1.017178083155432),(1.016324003281509),(1.015787884326783),(1.01607200764995),(1.017192416630256),(1.018496344113342))
$rng = $ws.Range("E2:E25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @((1.027159414871329),(1.027260784977316),(1.027325016012484),(1.027351691987983),(1.027356151822224),(1.027325373742885),(1.027193955042797),(1.026951973814643),(1.026783695275742),(1.026709186440694),(1.026681264296139),(1.02668726482421),(1.026706883407301),(1.026718938452259),(1.026788605615563),(1.026831866510928),(1.026856941281114),(1.026865464213467),(1.026827241424317),(1.026701113016217),(1.026620386682175),(1.026663316050686),(1.026829331791917),(1.027015760875279))
$rng = $ws.Range("I2:I25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @((1.023455495720206),(1.02398883736891),(1.024333773904531),(1.024478743047606),(1.024503081496555),(1.024335711156547),(1.023635775751023),(1.022401141955604),(1.02157727761383),(1.021220366198029),(1.021087768630182),(1.021116212360972),(1.021209406140363),(1.021266822655555),(1.021600961110152),(1.021810511921161),(1.021932722620425),(1.021974390434998),(1.021788030825683),(1.021181963569668),(1.020800762261449),(1.021002857429127),(1.021798189116898),(1.022720466188899))
$rng = $ws.Range("J2:J25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @((1.026765689217259),(1.027217014760895),(1.027508442857501),(1.027630811984352),(1.027651349627449),(1.027510078539949),(1.026918342273675),(1.025871014148735),(1.025169769486268),(1.024865420969711),(1.024752267727835),(1.02477654419784),(1.024856069815028),(1.024905054342714),(1.025189953349994),(1.025368475124646),(1.025472535753646),(1.025508006149387),(1.025349328462723),(1.024832654381688),(1.024507196019504),(1.024679784569091),(1.025357980225646),(1.026142312449795))
$rng = $ws.Range("K2:K25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @((1.022382269638418),(1.022947433595613),(1.023313247328722),(1.023467061657256),(1.023492889261529),(1.023315302501812),(1.022573244902177),(1.021266580030363),(1.020396173831303),(1.020019460912652),(1.019879561059593),(1.019909568761439),(1.020007896166747),(1.020068482686973),(1.020421178898886),(1.020642464642274),(1.020771554031709),(1.020815573050013),(1.020618720992351),(1.019978940410962),(1.019576848150657),(1.019789988981035),(1.020629449672835),(1.02160426565553))
$rng = $ws.Range("L2:L25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}

$vals = @((1.024908919916115),(1.025443018971135),(1.025788445356202),(1.025933620372088),(1.025957993384426),(1.025790385359338),(1.025089455965254),(1.023853068848231),(1.023028034524593),(1.022670616253643),(1.022537830382041),(1.022566314506191),(1.022659640631438),(1.022717138684667),(1.023051751654197),(1.023261600051392),(1.023383984303868),(1.023425711291517),(1.023239087030181),(1.022632159089145),(1.022250416431345),(1.022452798597341),(1.023249259747335),(1.024172846558564))
$rng = $ws.Range("N2:N25")
for ($i = 0; $i -lt $vals.Length; $i++) {
    $rng.Cells.Item($i + 1, 1).Value = $vals[$i]
}
